$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header F1
$ws.Range("F1").Value = "newColumn"

# Update the date for existing rows 2-6 and zero out column F
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 2).Value = 44146.03377314815
    $ws.Cells.Item($r, 6).Value = 0
}

# Source data for the new block of rows (same as rows 2-6 but new date/id pattern)
$ids     = @(1, 2, 3, 4, 5)
$regions = @("Alaska", "Illinois", "Toronto", "Pensilvania", "D.C. Columbia")
$brands  = @("ADIDAS", "NEW BALANCE", "TOMMY HILFIGER", "M TAC", "five eleven")
$states  = @("NY", "IL", "OH", "NV", "CA")

for ($i = 0; $i -lt 5; $i++) {
    $row = 7 + $i
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = 44146.03579861111
    $ws.Cells.Item($row, 3).Value = $regions[$i]
    $ws.Cells.Item($row, 4).Value = $brands[$i]
    $ws.Cells.Item($row, 5).Value = $states[$i]
    $ws.Cells.Item($row, 6).Value = 0
}

# Apply the same date style (style index 1, numFmtId 14) to the new B column cells
$ws.Range("B2:B6").Copy()
$ws.Range("B7:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
